# Add a new paragraph "Hello Zaid" right after the "Hello all" paragraph.
#
# In the original document, the second paragraph's text "Hello all" is
# immediately followed (still inside the same paragraph) by the hidden
# "_GoBack" bookmark that Word drops at the last edit location. The
# target edit turns that into two paragraphs:
#   Para 2: "Hello all"
#   Para 3: "Hello Zaid"   <- the "_GoBack" bookmark now trails this text
#
# Reproduce that the way a user typing at the end of "Hello all" would:
# insert the new sentence right at the existing bookmark (so the
# zero-width bookmark naturally slides to keep trailing the freshly
# typed text), then split the paragraph at the point just before that
# new text so it becomes its own paragraph.

$d = $word.ActiveDocument

# Remember where the bookmark currently sits (end of "Hello all", before
# the paragraph mark) so we can split the paragraph there afterwards.
$para2 = $d.Paragraphs(2)

$bm = $null
try {
    $bm = $d.Bookmarks("_GoBack")
} catch {
    $bm = $null
}

if ($bm -ne $null) {
    $insertPos = $bm.Range.Start
    # Typing directly into the bookmark's own (collapsed) range makes the
    # bookmark move/extend to keep trailing the inserted text.
    $bm.Range.InsertAfter("Hello Zaid") | Out-Null
} else {
    # Fallback: no bookmark present, just insert at the end of the
    # paragraph's text.
    $endRange = $para2.Range
    $endRange.MoveEnd(1, -1) | Out-Null   # wdCharacter: exclude the paragraph mark
    $insertPos = $endRange.End
    $endRange.Collapse(0) | Out-Null      # wdCollapseEnd
    $endRange.InsertAfter("Hello Zaid") | Out-Null
}

# Break the paragraph right before the text we just inserted so "Hello
# all" and "Hello Zaid" become separate paragraphs, with the bookmark
# staying attached to "Hello Zaid".
$splitPoint = $d.Range($insertPos, $insertPos)
$splitPoint.InsertParagraphAfter() | Out-Null
